$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column width changes
# Note: Excel's ColumnWidth (characters) is offset from the stored OOXML
# width by the default column padding (5/6 char for Calibri 11), so we
# subtract that offset to land on the exact target stored widths (12, 22).
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668

# Row 2
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 672.340305337043
$ws.Cells.Item(2,6).Value = 0

# Row 3
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 4992.1832
$ws.Cells.Item(3,6).Value = 0

# Row 4
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 142.502095025027
$ws.Cells.Item(4,6).Value = 0

# Row 6
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 106.82
$ws.Cells.Item(6,6).Value = 0

# Row 7
$ws.Cells.Item(7,3).Value = 2300
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 2300
$ws.Cells.Item(7,6).Value = 0

# Row 8
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 750
$ws.Cells.Item(8,6).Value = 0

# Row 10
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 650.25
$ws.Cells.Item(10,6).Value = 0

# Row 13
$ws.Cells.Item(13,3).Value = 130
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 130
$ws.Cells.Item(13,6).Value = 0

# Row 14
$ws.Cells.Item(14,3).Value = 240
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 240
$ws.Cells.Item(14,6).Value = 0

# Row 15
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 1505.12
$ws.Cells.Item(15,6).Value = 0

# Row 16
$ws.Cells.Item(16,3).Value = 40390.17
$ws.Cells.Item(16,4).Value = 793.77
$ws.Cells.Item(16,5).Value = 39596.4
$ws.Cells.Item(16,6).Value = 0.01965255407441959

# Row 17
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 342
$ws.Cells.Item(17,6).Value = 0

# Row 18
$ws.Cells.Item(18,3).Value = 2300
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 2300
$ws.Cells.Item(18,6).Value = 0

# Row 19 (TOTAL)
$ws.Cells.Item(19,4).Value = 793.77
$ws.Cells.Item(19,5).Value = 54615.93560036207
$ws.Cells.Item(19,6).Value = 0.01432546864127019
